$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting (bold, bordered, centered style) from G1:H1 to the
# new header cells I1:J1 before writing their values.
$ws.Range("G1:H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data columns I and J
$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 8
$ws.Range("J3").Value = 8

$ws.Range("I4").Value = 9
$ws.Range("J4").Value = 9

$ws.Range("I5").Value = 7
$ws.Range("J5").Value = 7

$ws.Range("I6").Value = 2
$ws.Range("J6").Value = 2

$ws.Range("I7").Value = 6
$ws.Range("J7").Value = 6

$ws.Range("I8").Value = 6
$ws.Range("J8").Value = 6
